$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "yyyyyyy"
$ws.Range("B3").Value = "123123WW!"
